# Commit: "medical panel added to week 10"
#
# The canonical-OOXML diff for this change shows the final slide of the
# deck -- "Case Study:  Assignment 1" (the Jimmy / 2-year-old case-study
# slide, sldId 297) -- being removed entirely: its part
# (ppt/slides/slide16.xml) disappears, and the corresponding <p:sldId>
# entry is dropped from <p:sldIdLst> in presentation.xml, leaving every
# other slide in the deck in the same order. Reproduce that with the
# standard PowerPoint Slide.Delete COM call.

$p = $ppt.ActivePresentation

# Find the "Case Study:  Assignment 1" slide by its title; fall back to
# the last slide in the deck (where it lives in before.pptx) if the
# title can't be matched for some reason.
$targetIndex = $p.Slides.Count
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.Shapes.HasTitle) {
        $title = $slide.Shapes.Title.TextFrame.TextRange.Text
        if ($title.Contains("Case Study") -and $title.Contains("Assignment 1")) {
            $targetIndex = $i
        }
    }
}

$p.Slides.Item($targetIndex).Delete()
